# Updated symbol list on Tue Dec 20 21:12:33 UTC 2022 with GitHub Actions
#
# Refreshes the crypto price/volume/hour snapshot in the sheet: new Price
# (column D) and Hora (column G) readings for the 21:xx UTC poll, plus a
# couple of "Worst/Best in 24h" label moves in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a target cell + its new literal value. D/G hold
# numeric-looking text ("250.15", "21", ...) that must remain text (the
# sheet stores them as inline strings), so those cells get NumberFormat
# "@" (Text) applied first -- otherwise assigning a numeric-looking string
# to .Value auto-converts it to a real number, same as interactive Excel.
$changes = @(
    @{ Cell = "D2"; Value = '250.15'; AsText = $true },
    @{ Cell = "G2"; Value = '21'; AsText = $true },
    @{ Cell = "D3"; Value = '23.14'; AsText = $true },
    @{ Cell = "G3"; Value = '21'; AsText = $true },
    @{ Cell = "D4"; Value = '5.450'; AsText = $true },
    @{ Cell = "G4"; Value = '21'; AsText = $true },
    @{ Cell = "D5"; Value = '0.05665'; AsText = $true },
    @{ Cell = "G5"; Value = '21'; AsText = $true },
    @{ Cell = "D6"; Value = '3.414'; AsText = $true },
    @{ Cell = "G6"; Value = '21'; AsText = $true },
    @{ Cell = "D7"; Value = '6.395'; AsText = $true },
    @{ Cell = "G7"; Value = '21'; AsText = $true },
    @{ Cell = "D8"; Value = '0.8148'; AsText = $true },
    @{ Cell = "G8"; Value = '21'; AsText = $true },
    @{ Cell = "D9"; Value = '0.9334'; AsText = $true },
    @{ Cell = "G9"; Value = '21'; AsText = $true },
    @{ Cell = "G10"; Value = '21'; AsText = $true },
    @{ Cell = "D11"; Value = '0.07535'; AsText = $true },
    @{ Cell = "G11"; Value = '21'; AsText = $true },
    @{ Cell = "D12"; Value = '0.03122'; AsText = $true },
    @{ Cell = "G12"; Value = '21'; AsText = $true },
    @{ Cell = "D13"; Value = '0.03096'; AsText = $true },
    @{ Cell = "G13"; Value = '21'; AsText = $true },
    @{ Cell = "D14"; Value = '0.09362'; AsText = $true },
    @{ Cell = "G14"; Value = '21'; AsText = $true },
    @{ Cell = "D15"; Value = '3.758'; AsText = $true },
    @{ Cell = "G15"; Value = '21'; AsText = $true },
    @{ Cell = "D16"; Value = '0.001585'; AsText = $true },
    @{ Cell = "G16"; Value = '21'; AsText = $true },
    @{ Cell = "D17"; Value = '0.04755'; AsText = $true },
    @{ Cell = "G17"; Value = '21'; AsText = $true },
    @{ Cell = "D18"; Value = '0.0005794'; AsText = $true },
    @{ Cell = "E18"; Value = '17OneONE'; AsText = $false },
    @{ Cell = "G18"; Value = '21'; AsText = $true },
    @{ Cell = "D19"; Value = '0.006377'; AsText = $true },
    @{ Cell = "G19"; Value = '21'; AsText = $true },
    @{ Cell = "D20"; Value = '0.005042'; AsText = $true },
    @{ Cell = "G20"; Value = '21'; AsText = $true },
    @{ Cell = "D21"; Value = '0.001031'; AsText = $true },
    @{ Cell = "G21"; Value = '21'; AsText = $true },
    @{ Cell = "D22"; Value = '0.0001501'; AsText = $true },
    @{ Cell = "G22"; Value = '21'; AsText = $true },
    @{ Cell = "G23"; Value = '21'; AsText = $true },
    @{ Cell = "D24"; Value = '2.170'; AsText = $true },
    @{ Cell = "G24"; Value = '21'; AsText = $true },
    @{ Cell = "D25"; Value = '0.3301'; AsText = $true },
    @{ Cell = "G25"; Value = '21'; AsText = $true },
    @{ Cell = "D26"; Value = '0.1266'; AsText = $true },
    @{ Cell = "G26"; Value = '21'; AsText = $true },
    @{ Cell = "G27"; Value = '21'; AsText = $true },
    @{ Cell = "D28"; Value = '0.0003002'; AsText = $true },
    @{ Cell = "G28"; Value = '21'; AsText = $true },
    @{ Cell = "G29"; Value = '21'; AsText = $true },
    @{ Cell = "G30"; Value = '21'; AsText = $true },
    @{ Cell = "G31"; Value = '21'; AsText = $true },
    @{ Cell = "G32"; Value = '21'; AsText = $true },
    @{ Cell = "G33"; Value = '21'; AsText = $true },
    @{ Cell = "G34"; Value = '21'; AsText = $true },
    @{ Cell = "G35"; Value = '21'; AsText = $true },
    @{ Cell = "G36"; Value = '21'; AsText = $true },
    @{ Cell = "G37"; Value = '21'; AsText = $true },
    @{ Cell = "G38"; Value = '21'; AsText = $true },
    @{ Cell = "G39"; Value = '21'; AsText = $true },
    @{ Cell = "D40"; Value = '0.04024'; AsText = $true },
    @{ Cell = "G40"; Value = '21'; AsText = $true },
    @{ Cell = "D41"; Value = '0.006776'; AsText = $true },
    @{ Cell = "G41"; Value = '21'; AsText = $true },
    @{ Cell = "D42"; Value = '0.1070'; AsText = $true },
    @{ Cell = "G42"; Value = '21'; AsText = $true },
    @{ Cell = "D43"; Value = '0.002712'; AsText = $true },
    @{ Cell = "G43"; Value = '21'; AsText = $true },
    @{ Cell = "D44"; Value = '0.008026'; AsText = $true },
    @{ Cell = "G44"; Value = '21'; AsText = $true },
    @{ Cell = "D45"; Value = '0.00005807'; AsText = $true },
    @{ Cell = "G45"; Value = '21'; AsText = $true },
    @{ Cell = "G46"; Value = '21'; AsText = $true },
    @{ Cell = "D47"; Value = '0.5003'; AsText = $true },
    @{ Cell = "E47"; Value = '46CoinbaseStockTokenCOINWorstin24h'; AsText = $false },
    @{ Cell = "G47"; Value = '21'; AsText = $true },
    @{ Cell = "G48"; Value = '21'; AsText = $true },
    @{ Cell = "D49"; Value = '0.00002102'; AsText = $true },
    @{ Cell = "G49"; Value = '21'; AsText = $true },
    @{ Cell = "G50"; Value = '21'; AsText = $true },
    @{ Cell = "G51"; Value = '21'; AsText = $true }
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    if ($chg.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $chg.Value
}
